$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column E (Meas act SPECT) to a custom width (stored width == 21)
$ws.Columns.Item(5).ColumnWidth = 20.17

# Update the recalculated "Meas act SPECT" values for rows 17-21
$ws.Range("E17").Value = 0.0654805076123525
$ws.Range("E18").Value = 0.1260997997464052
$ws.Range("E19").Value = 0.25816899703943336
$ws.Range("E20").Value = 0.5030108691473354
$ws.Range("E21").Value = 3.4604660149691666

# Move the active selection to K20, matching the author's last cursor position
$ws.Range("K20").Select()
